$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental row (B7): was blank, now holds the literal text "true".
# A direct .Value = "true" assignment gets auto-coerced to a Boolean by
# the engine (same as typing into Excel), so instead build it as a text
# formula and collapse it to a literal value via copy / paste-special —
# that keeps the cell's type as text (t="s") instead of boolean (t="b").
$b7 = $ws.Range("B7")
$b7.Formula = "=""true"""
$b7.Copy()
$b7.PasteSpecial(-4163)

# Date row (B8): refresh the generated timestamp.
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
